$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Päävuokraukset")

$ws.Range("J3").Value = "test2@test.test / test3@test.test"
$ws.Range("J4").Value = "asd"
